# Update countries & provincias Spain
# Applies the COVID-19 "paises" data refresh:
#  - updates the "datos actualizados" timestamp
#  - refreshes Casos totales / Nuevos casos / Casos activos / Recuperados /
#    Casos criticos / Muertes hoy / Muertes for the countries whose figures
#    moved, including the handful of countries that swapped rank (and thus
#    row) with a neighbour because of the new totals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 11 de Agosto de 2020 a las 00:07"

# --- Per-row data refresh ----------------------------------------------
# Row = worksheet row number; Name/B..H = País / Casos totales / Nuevos
# casos / Casos activos / Recuperados / Casos críticos / Muertes hoy / Muertes
$rows = @(
  @{ Row=4;   Name="Estados Unidos";   B=5239825; C=40381; D=2682423; E=2391349; F=0; G=436; H=166053 },
  @{ Row=5;   Name="Brasil";           B=3057470; C=21888; D=2118460; E=837258;  F=0; G=616; H=101752 },
  @{ Row=6;   Name="India";            B=2267153; C=53016; D=1581640; E=640160;  F=0; G=887; H=45353  },
  @{ Row=8;   Name="Sudafrica";        B=563598;  C=3739;  D=417200;  E=135777;  F=0; G=213; H=10621  },
  @{ Row=30;  Name="Egipto";           B=95666;   C=174;   D=53779;   E=36852;   F=0; G=26;  H=5035   },
  @{ Row=46;  Name="Guatemala";        B=56987;   C=382;   D=45589;   E=9176;    F=0; G=11;  H=2222   },
  @{ Row=58;  Name="Suiza";            B=36708;   C=105;   D=32400;   E=2321;    F=0; G=1;   H=1987   },
  @{ Row=62;  Name="Uzbekistan";       B=31304;   C=695;   D=22992;   E=8112;    F=0; G=6;   H=200    },
  @{ Row=76;  Name="Costa de Marfil";  B=16798;   C=83;    D=13052;   E=3641;    F=0; G=0;   H=105    },
  @{ Row=81;  Name="Bulgaria";         B=13512;   C=116;   D=7980;    E=5073;    F=0; G=12;  H=459    },
  @{ Row=91;  Name="Gabon";            B=8006;    C=83;    D=5823;    E=2132;    F=0; G=0;   H=51     },
  @{ Row=92;  Name="Guinea";           B=7930;    C=0;     D=6898;    E=982;     F=0; G=0;   H=50     },
  @{ Row=121; Name="Cabo Verde";       B=2883;    C=25;    D=2128;    E=723;     F=0; G=0;   H=32     },
  @{ Row=122; Name="Sri Lanka";        B=2871;    C=27;    D=2593;    E=267;     F=0; G=0;   H=11     },
  @{ Row=136; Name="Yemen";            B=1832;    C=28;    D=915;     E=399;     F=0; G=3;   H=518    },
  @{ Row=152; Name="Togo";             B=1067;    C=7;     D=729;     E=313;     F=0; G=2;   H=25     },
  @{ Row=202; Name="Timor Oriental";   B=25;      C=0;     D=24;      E=1;       F=0; G=0;   H=0      },
  @{ Row=203; Name="Santa Lucia";      B=25;      C=0;     D=24;      E=1;       F=0; G=0;   H=0      },
  @{ Row=213; Name="Montserrat";       B=13;      C=0;     D=12;      E=0;       F=0; G=0;   H=1      },
  @{ Row=214; Name="Islas Malvinas";   B=13;      C=0;     D=13;      E=0;       F=0; G=0;   H=0      }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Name
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
    $ws.Cells.Item($r.Row, 7).Value = $r.G
    $ws.Cells.Item($r.Row, 8).Value = $r.H
}
